$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FIXE")

# Insert a new row at 83 (shifts old rows 83-90 down to 84-91, carrying
# their values/formats with them, same as Excel's native row insert).
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new "ADNP75" field
# (write the libelle text before the short code so the shared-strings
# table picks up "Non Programmé (NP)" ahead of "ADNP75", matching the
# order they were appended to the workbook).
$ws.Range("F83").Value = "Non Programmé (NP)"
$ws.Range("A83").Value = 206
$ws.Range("B83").Value = "ADNP75"
$ws.Range("C83").Value = 1
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = "c"

# The row that used to be row 83 (FILLER5) is now row 84: it shrinks from
# 4 bytes to 3 since the new ADNP75 field took 1 of its bytes, and its
# starting "position" shifts from 206 to 207.
$ws.Range("A84").Value = 207
$ws.Range("C84").Value = 3

# Match the author's final cursor position/selection.
$ws.Activate() | Out-Null
$ws.Range("A84").Select() | Out-Null
